$p = $ppt.ActivePresentation

# --- Slide 11: "Программный код. Очистка заполненных линий" ---
# Merge the two runs "О" + "чистка " into a single run "Очистка "
$s11 = $p.Slides.Item(11)
$tr11 = $s11.Shapes.Item(1).TextFrame.TextRange
$merge11 = $tr11.Characters(18, 8)
$merge11.Text = "Очистка "

# --- Slide 12: "Программный код. Отрисовка поля" ---
# Merge the two runs "О" + "трисовка" into a single run "Отрисовка"
$s12 = $p.Slides.Item(12)
$tr12 = $s12.Shapes.Item(1).TextFrame.TextRange
$merge12 = $tr12.Characters(18, 9)
$merge12.Text = "Отрисовка"

# --- Slide 2: "Структура" ---
# The second shape (content placeholder) was empty; add two paragraphs of text.
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange

$run1 = $tr2.InsertAfter("2 класса")
$run1.LanguageID = "ru-RU"

$run2 = $tr2.InsertAfter([char]13 + "12 функций")
$run2.LanguageID = "ru-RU"
